# Applies the 6.4.2.1 workbook update:
#  - Retitle the header text in B1 (drop the period after "6.4.2.1")
#  - Update several data values in column L (2022)
#  - Replace the formula in L7 with its literal resulting value
#  - Move the active selection to O2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header label: "6.4.2.1. Общий объем забора пресной воды " -> "6.4.2.1 Общий объем забора пресной воды "
$ws.Range("B1").Value = "6.4.2.1 Общий объем забора пресной воды "

# Column L (year 2022) value updates
$ws.Range("L5").Value = 8741.9
$ws.Range("L7").Value = 8483.5
$ws.Range("L14").Value = 1327.6
$ws.Range("L18").Value = 54

# Move selection to O2
$ws.Range("O2").Select() | Out-Null
